$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.291993498802185
$ws.Range("B1").Value = 1.668675899505615
$ws.Range("C1").Value = 2.399070978164673
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.321151971817017
